$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 68: Sno 67, date 21-Jul-2022, Code, segformer attention map ---
$ws.Range("A67:G67").Copy()
$ws.Range("A68:G68").PasteSpecial(-4122)

$ws.Range("A68").Value = 67
$ws.Range("B68").Value = 44763
$ws.Range("C68").Value = 0.34722222222222227
$ws.Range("D68").Value = 0.375
$ws.Range("E68").Formula = "=D68-C68"
$ws.Range("F68").Value = "Code"
$ws.Range("G68").Value = "1. segformer starter nb attention map ppt and sample code"

# --- New row 69: Sno 68, date 22-Jul-2022, Literature survey, visualize attention map ---
$ws.Range("A67:G67").Copy()
$ws.Range("A69:G69").PasteSpecial(-4122)

$ws.Range("A69").Value = 68
$ws.Range("B69").Value = 44764
$ws.Range("C69").Value = 0.34722222222222227
$ws.Range("D69").Value = 0.375
$ws.Range("E69").Formula = "=D69-C69"
$ws.Range("F69").Value = "Literature survey"
$ws.Range("G69").Value = "1. visualize attention map and models output literature survey"

# --- Move "Total Hours" row from 69 down to row 75 (now that 69 holds real entries) ---
$ws.Range("C75").Value = "Total Hours"
$ws.Range("C75").HorizontalAlignment = -4108
$ws.Range("C75").VerticalAlignment = -4108

$ws.Range("E75").Formula = "=SUM(E2:E74)"
$ws.Range("E75").NumberFormat = "[hh]:mm"
$ws.Range("E75").HorizontalAlignment = -4108
$ws.Range("E75").VerticalAlignment = -4108

# --- Update the view to reflect scrolling to the bottom of the (now longer) sheet ---
$ws.Application.ActiveWindow.ScrollRow = 61
$ws.Range("E76").Select()

Write-Host "edit applied"
